$d = $word.ActiveDocument

# --- Remove the two paragraphs that are dropped entirely -------------------
# (original paragraph 4: "ADMIN: Can enroll more students than a section
#  capacity allows." and original paragraph 3: "ADMIN: Students created via
#  browser cannot be assigned to sections." -- delete from the bottom up so
#  indices of paragraphs above stay valid)
$d.Paragraphs.Item(4).Range.Delete()
$d.Paragraphs.Item(3).Range.Delete()

# --- Paragraph 2 ("ADMIN: Can't add student classification...") becomes
#     "ADMIN: Students created via browser cannot be assigned to sections."
$d.Paragraphs.Item(2).Range.Find.Execute(
    "ADMIN: Can’t add student classification when creating students.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ADMIN: Students created via browser cannot be assigned to sections.", 2)

# --- Move the _GoBack bookmark to the start of that paragraph (a document
#     may only have one bookmark named "_GoBack"; adding it here relocates
#     it away from its old spot in the last paragraph automatically).
$bookmarkRange = $d.Paragraphs.Item(2).Range.Duplicate
$bookmarkRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- Last paragraph: split "...make MultipleChoice and Matching..." so the
#     text reads "...make MultipleChoice and Matching..." still, but broken
#     into three runs around "MultipleChoice" (prefix / word / suffix).
$d.Paragraphs.Item(3).Range.Find.Execute(
    "INSTRUCTOR/STUDENT: PM said to make MultipleChoice and Matching question point value to show total for question not just for each correct choice.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "INSTRUCTOR/STUDENT: PM said to make MultipleChoice and Matching question point value to show total for question not just for each correct choice.", 2)
